$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update header row text (remove the space in "Data N" labels; "Total" unchanged)
$ws.Range("A1").Value = "Data1"
$ws.Range("B1").Value = "Data2"
$ws.Range("C1").Value = "Data3"
$ws.Range("D1").Value = "Total"

# Move the active selection to D1 (was D9)
$ws.Range("D1").Select() | Out-Null
